# Apply updated crypto price/volume figures (Fri Aug 30 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.834.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.510.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.87%  '

$ws.Range("E4").Value = '  +0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.96%  '

$ws.Range("E7").Value = '  +0.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.562'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.512.06'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1000'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.58%  '

$ws.Range("E11").Value = '  +0.71%  '

$ws.Range("E12").Value = '  -2.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.954.71'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.15%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.63%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.803.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000139'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.51%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.508.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.58%  '

$ws.Range("E19").Value = '  -1.74%  '

$ws.Range("E20").Value = '  -0.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.73'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.87%  '

$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("E23").Value = '  +0.47%  '

$ws.Range("E24").Value = '  +0.34%  '

$ws.Range("E25").Value = '  -3.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.166'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.67'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0766'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.38%  '

$ws.Range("E31").Value = '  -1.61%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '163.63'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.06%  '

$ws.Range("E33").Value = '  +0.22%  '

$ws.Range("E34").Value = '  -5.73%  '

$ws.Range("E35").Value = '  -3.36%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.39'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.21'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.47%  '

$ws.Range("E38").Value = '  -3.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.79'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.58%  '

$ws.Range("E40").Value = '  -1.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.799'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.98%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.17'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '277.96'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.88%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.594'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.10%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0932'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '121.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0509'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0222'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.44%  '
